$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12; this shifts existing rows 12-50 down to 13-51
$ws.Rows.Item(12).Insert()

# Fill the new row 12 with the latest weekly price observation for Orégano
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44608
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = 100112029
$ws.Cells.Item(12, 7).Value = "Orégano"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 16
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = 9500
$ws.Cells.Item(12, 14).Value = "$/docena de atados"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 3167
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = "Hortaliza"
